# Update symbol list cell values (prices, volumes, coin metadata)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'248.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("D4").Value = "'5.397"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05681"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.397"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.312"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8050"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9216"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1409"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07421"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03089"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03019"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09368"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Value = "'0.001577"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04734"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.01827"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.0005849"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18OneONEWorstin24h"
$ws.Range("D20").Value = "'0.006467"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.004983"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Value = "'0.0001499"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Value = "'2.175"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Value = "'0.1298"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.03995"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006787"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1067"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002715"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.007518"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005800"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.3999"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.2136"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.01010"
$ws.Range("D50").Style = "Normal"
